$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 1.02
$ws.Cells.Item(2, 3).Value2 = 1.054538976934235
$ws.Cells.Item(2, 4).Value2 = 1.059391456928356
$ws.Cells.Item(2, 5).Value2 = 1.058222750956105
$ws.Cells.Item(2, 6).Value2 = 1.068813477642035
$ws.Cells.Item(2, 9).Value2 = 1.052610208696019
$ws.Cells.Item(2, 10).Value2 = 1.059550545316738
$ws.Cells.Item(2, 11).Value2 = 1.062121088271499
$ws.Cells.Item(2, 12).Value2 = 1.060955573782085
$ws.Cells.Item(2, 13).Value2 = 1.071517655283849
$ws.Cells.Item(2, 14).Value2 = 1.023510039378259
$ws.Cells.Item(3, 2).Value2 = 1.02
$ws.Cells.Item(3, 3).Value2 = 1.055668169955069
$ws.Cells.Item(3, 4).Value2 = 1.060286546329833
$ws.Cells.Item(3, 5).Value2 = 1.059302154456931
$ws.Cells.Item(3, 6).Value2 = 1.069850032242391
$ws.Cells.Item(3, 9).Value2 = 1.053000665611775
$ws.Cells.Item(3, 10).Value2 = 1.060330189471074
$ws.Cells.Item(3, 11).Value2 = 1.062830193318865
$ws.Cells.Item(3, 12).Value2 = 1.061848294857886
$ws.Cells.Item(3, 13).Value2 = 1.072369712587539
$ws.Cells.Item(3, 14).Value2 = 1.023781149124378
$ws.Cells.Item(4, 2).Value2 = 1.02
$ws.Cells.Item(4, 3).Value2 = 1.056398781197901
$ws.Cells.Item(4, 4).Value2 = 1.060865660457973
$ws.Cells.Item(4, 5).Value2 = 1.060000877362404
$ws.Cells.Item(4, 6).Value2 = 1.070520996930373
$ws.Cells.Item(4, 9).Value2 = 1.053252026146662
$ws.Cells.Item(4, 10).Value2 = 1.060834047281768
$ws.Cells.Item(4, 11).Value2 = 1.063288327191
$ws.Cells.Item(4, 12).Value2 = 1.06242562404714
$ws.Cells.Item(4, 13).Value2 = 1.072920685891417
$ws.Cells.Item(4, 14).Value2 = 1.023956102967869
$ws.Cells.Item(5, 2).Value2 = 1.02
$ws.Cells.Item(5, 3).Value2 = 1.056705918564615
$ws.Cells.Item(5, 4).Value2 = 1.06110910366836
$ws.Cells.Item(5, 5).Value2 = 1.060294687079864
$ws.Cells.Item(5, 6).Value2 = 1.070803129189773
$ws.Cells.Item(5, 9).Value2 = 1.053357389172399
$ws.Cells.Item(5, 10).Value2 = 1.061045720028916
$ws.Cells.Item(5, 11).Value2 = 1.063480758210507
$ws.Cells.Item(5, 12).Value2 = 1.062668256731436
$ws.Cells.Item(5, 13).Value2 = 1.073152227944392
$ws.Cells.Item(5, 14).Value2 = 1.02402954047706
$ws.Cells.Item(6, 2).Value2 = 1.02
$ws.Cells.Item(6, 3).Value2 = 1.056757487625551
$ws.Cells.Item(6, 4).Value2 = 1.061149977896052
$ws.Cells.Item(6, 5).Value2 = 1.060344022941113
$ws.Cells.Item(6, 6).Value2 = 1.070850503884692
$ws.Cells.Item(6, 9).Value2 = 1.053375061976519
$ws.Cells.Item(6, 10).Value2 = 1.061081252082348
$ws.Cells.Item(6, 11).Value2 = 1.063513058352693
$ws.Cells.Item(6, 12).Value2 = 1.062708991350619
$ws.Cells.Item(6, 13).Value2 = 1.073191099770697
$ws.Cells.Item(6, 14).Value2 = 1.024041864325997
$ws.Cells.Item(7, 2).Value2 = 1.02
$ws.Cells.Item(7, 3).Value2 = 1.056402885226887
$ws.Cells.Item(7, 4).Value2 = 1.060868913422325
$ws.Cells.Item(7, 5).Value2 = 1.06000480300003
$ws.Cells.Item(7, 6).Value2 = 1.07052476656456
$ws.Cells.Item(7, 9).Value2 = 1.053253435226285
$ws.Cells.Item(7, 10).Value2 = 1.060836876248725
$ws.Cells.Item(7, 11).Value2 = 1.063290899125382
$ws.Cells.Item(7, 12).Value2 = 1.062428866417119
$ws.Cells.Item(7, 13).Value2 = 1.07292378010975
$ws.Cells.Item(7, 14).Value2 = 1.023957084687835
$ws.Cells.Item(8, 2).Value2 = 1.02
$ws.Cells.Item(8, 3).Value2 = 1.054920604331892
$ws.Cells.Item(8, 4).Value2 = 1.059693970807909
$ws.Cells.Item(8, 5).Value2 = 1.05858748336328
$ws.Cells.Item(8, 6).Value2 = 1.069163735582717
$ws.Cells.Item(8, 9).Value2 = 1.052742432893583
$ws.Cells.Item(8, 10).Value2 = 1.059814159130742
$ws.Cells.Item(8, 11).Value2 = 1.062360879636596
$ws.Cells.Item(8, 12).Value2 = 1.06125733973425
$ws.Cells.Item(8, 13).Value2 = 1.071805687681002
$ws.Cells.Item(8, 14).Value2 = 1.023601759965686
$ws.Cells.Item(9, 2).Value2 = 1.02
$ws.Cells.Item(9, 3).Value2 = 1.052308196757263
$ws.Cells.Item(9, 4).Value2 = 1.057623044054562
$ws.Cells.Item(9, 5).Value2 = 1.05609208678919
$ws.Cells.Item(9, 6).Value2 = 1.066767290414286
$ws.Cells.Item(9, 9).Value2 = 1.051832084516206
$ws.Cells.Item(9, 10).Value2 = 1.058007205893485
$ws.Cells.Item(9, 11).Value2 = 1.060716661394969
$ws.Cells.Item(9, 12).Value2 = 1.059190485140425
$ws.Cells.Item(9, 13).Value2 = 1.069832663169135
$ws.Cells.Item(9, 14).Value2 = 1.022972015894957
$ws.Cells.Item(10, 2).Value2 = 1.02
$ws.Cells.Item(10, 3).Value2 = 1.050566238105772
$ws.Cells.Item(10, 4).Value2 = 1.056242064880391
$ws.Cells.Item(10, 5).Value2 = 1.054429875225399
$ws.Cells.Item(10, 6).Value2 = 1.065170900730469
$ws.Cells.Item(10, 9).Value2 = 1.051218518819687
$ws.Cells.Item(10, 10).Value2 = 1.056799322892147
$ws.Cells.Item(10, 11).Value2 = 1.05961686641723
$ws.Cells.Item(10, 12).Value2 = 1.057810893026301
$ws.Cells.Item(10, 13).Value2 = 1.068515413453977
$ws.Cells.Item(10, 14).Value2 = 1.022549754858937
$ws.Cells.Item(11, 2).Value2 = 1.02
$ws.Cells.Item(11, 3).Value2 = 1.049811853130659
$ws.Cells.Item(11, 4).Value2 = 1.055643997594738
$ws.Cells.Item(11, 5).Value2 = 1.053710441146936
$ws.Cells.Item(11, 6).Value2 = 1.064479936585767
$ws.Cells.Item(11, 9).Value2 = 1.05095125352985
$ws.Cells.Item(11, 10).Value2 = 1.056275519949623
$ws.Cells.Item(11, 11).Value2 = 1.059139773274333
$ws.Cells.Item(11, 12).Value2 = 1.057213107399825
$ws.Cells.Item(11, 13).Value2 = 1.06794457481453
$ws.Cells.Item(11, 14).Value2 = 1.022366334047013
$ws.Cells.Item(12, 2).Value2 = 1.02
$ws.Cells.Item(12, 3).Value2 = 1.049531623838671
$ws.Cells.Item(12, 4).Value2 = 1.055421834317477
$ws.Cells.Item(12, 5).Value2 = 1.053443257847335
$ws.Cells.Item(12, 6).Value2 = 1.064223323820082
$ws.Cells.Item(12, 9).Value2 = 1.050851740470813
$ws.Cells.Item(12, 10).Value2 = 1.056080838169789
$ws.Cells.Item(12, 11).Value2 = 1.058962427777214
$ws.Cells.Item(12, 12).Value2 = 1.056991000481457
$ws.Cells.Item(12, 13).Value2 = 1.067732470221105
$ws.Cells.Item(12, 14).Value2 = 1.022298116327154
$ws.Cells.Item(13, 2).Value2 = 1.02
$ws.Cells.Item(13, 3).Value2 = 1.049591734768633
$ws.Cells.Item(13, 4).Value2 = 1.055469489753824
$ws.Cells.Item(13, 5).Value2 = 1.05350056747683
$ws.Cells.Item(13, 6).Value2 = 1.064278366232362
$ws.Cells.Item(13, 9).Value2 = 1.050873097173086
$ws.Cells.Item(13, 10).Value2 = 1.056122603424038
$ws.Cells.Item(13, 11).Value2 = 1.059000474962641
$ws.Cells.Item(13, 12).Value2 = 1.057038646018084
$ws.Cells.Item(13, 13).Value2 = 1.06777797053835
$ws.Cells.Item(13, 14).Value2 = 1.022312753200658
$ws.Cells.Item(14, 2).Value2 = 1.02
$ws.Cells.Item(14, 3).Value2 = 1.04978868964673
$ws.Cells.Item(14, 4).Value2 = 1.055625633805664
$ws.Cells.Item(14, 5).Value2 = 1.053688354746484
$ws.Cells.Item(14, 6).Value2 = 1.06445872403865
$ws.Cells.Item(14, 9).Value2 = 1.050943032623476
$ws.Cells.Item(14, 10).Value2 = 1.05625942990459
$ws.Cells.Item(14, 11).Value2 = 1.059125116532769
$ws.Cells.Item(14, 12).Value2 = 1.057194749255904
$ws.Cells.Item(14, 13).Value2 = 1.067927043607391
$ws.Cells.Item(14, 14).Value2 = 1.022360696925515
$ws.Cells.Item(15, 2).Value2 = 1.02
$ws.Cells.Item(15, 3).Value2 = 1.049910037739897
$ws.Cells.Item(15, 4).Value2 = 1.055721837384318
$ws.Cells.Item(15, 5).Value2 = 1.053804062828479
$ws.Cells.Item(15, 6).Value2 = 1.064569853998489
$ws.Cells.Item(15, 9).Value2 = 1.050986090504615
$ws.Cells.Item(15, 10).Value2 = 1.056343717554712
$ws.Cells.Item(15, 11).Value2 = 1.059201894831363
$ws.Cells.Item(15, 12).Value2 = 1.057290921282273
$ws.Cells.Item(15, 13).Value2 = 1.068018883194373
$ws.Cells.Item(15, 14).Value2 = 1.022390225097724
$ws.Cells.Item(16, 2).Value2 = 1.02
$ws.Cells.Item(16, 3).Value2 = 1.050616301821552
$ws.Cells.Item(16, 4).Value2 = 1.056281754625757
$ws.Cells.Item(16, 5).Value2 = 1.054477628262156
$ws.Cells.Item(16, 6).Value2 = 1.06521676367152
$ws.Cells.Item(16, 9).Value2 = 1.051236222863752
$ws.Cells.Item(16, 10).Value2 = 1.056834069447887
$ws.Cells.Item(16, 11).Value2 = 1.059648511051058
$ws.Cells.Item(16, 12).Value2 = 1.057850557329659
$ws.Cells.Item(16, 13).Value2 = 1.068553288385127
$ws.Cells.Item(16, 14).Value2 = 1.022561915688873
$ws.Cells.Item(17, 2).Value2 = 1.02
$ws.Cells.Item(17, 3).Value2 = 1.051059293495322
$ws.Cells.Item(17, 4).Value2 = 1.056632950484666
$ws.Cells.Item(17, 5).Value2 = 1.054900221769541
$ws.Cells.Item(17, 6).Value2 = 1.065622628498845
$ws.Cells.Item(17, 9).Value2 = 1.051392699067911
$ws.Cells.Item(17, 10).Value2 = 1.057141444537579
$ws.Cells.Item(17, 11).Value2 = 1.059928427114923
$ws.Cells.Item(17, 12).Value2 = 1.058201491168752
$ws.Cells.Item(17, 13).Value2 = 1.068888382934089
$ws.Cells.Item(17, 14).Value2 = 1.022669457593338
$ws.Cells.Item(18, 2).Value2 = 1.02
$ws.Cells.Item(18, 3).Value2 = 1.051317672962292
$ws.Cells.Item(18, 4).Value2 = 1.056837788020125
$ws.Cells.Item(18, 5).Value2 = 1.055146743804102
$ws.Cells.Item(18, 6).Value2 = 1.065859389933651
$ws.Cells.Item(18, 9).Value2 = 1.051483815863444
$ws.Cells.Item(18, 10).Value2 = 1.057320655851974
$ws.Cells.Item(18, 11).Value2 = 1.060091613057584
$ws.Cells.Item(18, 12).Value2 = 1.058406145121692
$ws.Cells.Item(18, 13).Value2 = 1.06908379348443
$ws.Cells.Item(18, 14).Value2 = 1.022732129107034
$ws.Cells.Item(19, 2).Value2 = 1.02
$ws.Cells.Item(19, 3).Value2 = 1.051405772012679
$ws.Cells.Item(19, 4).Value2 = 1.056907630797601
$ws.Cells.Item(19, 5).Value2 = 1.055230806623377
$ws.Cells.Item(19, 6).Value2 = 1.065940124116874
$ws.Cells.Item(19, 9).Value2 = 1.05151485837395
$ws.Cells.Item(19, 10).Value2 = 1.057381749512692
$ws.Cells.Item(19, 11).Value2 = 1.060147240951536
$ws.Cells.Item(19, 12).Value2 = 1.058475920096883
$ws.Cells.Item(19, 13).Value2 = 1.069150415910401
$ws.Cells.Item(19, 14).Value2 = 1.022753489001155
$ws.Cells.Item(20, 2).Value2 = 1.02
$ws.Cells.Item(20, 3).Value2 = 1.051011765733056
$ws.Cells.Item(20, 4).Value2 = 1.056595271432649
$ws.Cells.Item(20, 5).Value2 = 1.054854878346947
$ws.Cells.Item(20, 6).Value2 = 1.065579080218361
$ws.Cells.Item(20, 9).Value2 = 1.05137592649322
$ws.Cells.Item(20, 10).Value2 = 1.05710847390011
$ws.Cells.Item(20, 11).Value2 = 1.059898403499436
$ws.Cells.Item(20, 12).Value2 = 1.058163843412146
$ws.Cells.Item(20, 13).Value2 = 1.068852435059246
$ws.Cells.Item(20, 14).Value2 = 1.022657925143091
$ws.Cells.Item(21, 2).Value2 = 1.02
$ws.Cells.Item(21, 3).Value2 = 1.049730691827142
$ws.Cells.Item(21, 4).Value2 = 1.055579653675783
$ws.Cells.Item(21, 5).Value2 = 1.053633054796659
$ws.Cells.Item(21, 6).Value2 = 1.064405612001948
$ws.Cells.Item(21, 9).Value2 = 1.050922444972183
$ws.Cells.Item(21, 10).Value2 = 1.05621914117991
$ws.Cells.Item(21, 11).Value2 = 1.059088416340633
$ws.Cells.Item(21, 12).Value2 = 1.057148782479751
$ws.Cells.Item(21, 13).Value2 = 1.067883147232137
$ws.Cells.Item(21, 14).Value2 = 1.022346581097497
$ws.Cells.Item(22, 2).Value2 = 1.02
$ws.Cells.Item(22, 3).Value2 = 1.048925130109521
$ws.Cells.Item(22, 4).Value2 = 1.054941010677446
$ws.Cells.Item(22, 5).Value2 = 1.052865114813958
$ws.Cells.Item(22, 6).Value2 = 1.063668049213178
$ws.Cells.Item(22, 9).Value2 = 1.050635941166357
$ws.Cells.Item(22, 10).Value2 = 1.055659299128173
$ws.Cells.Item(22, 11).Value2 = 1.058578382132699
$ws.Cells.Item(22, 12).Value2 = 1.056510210181752
$ws.Cells.Item(22, 13).Value2 = 1.067273313971952
$ws.Cells.Item(22, 14).Value2 = 1.022150322978913
$ws.Cells.Item(23, 2).Value2 = 1.02
$ws.Cells.Item(23, 3).Value2 = 1.049352183499447
$ws.Cells.Item(23, 4).Value2 = 1.055279575461104
$ws.Cells.Item(23, 5).Value2 = 1.053272189066442
$ws.Cells.Item(23, 6).Value2 = 1.064059022296533
$ws.Cells.Item(23, 9).Value2 = 1.050787953368925
$ws.Cells.Item(23, 10).Value2 = 1.055956146963436
$ws.Cells.Item(23, 11).Value2 = 1.058848833404094
$ws.Cells.Item(23, 12).Value2 = 1.056848764083843
$ws.Cells.Item(23, 13).Value2 = 1.067596636529847
$ws.Cells.Item(23, 14).Value2 = 1.022254410896793
$ws.Cells.Item(24, 2).Value2 = 1.02
$ws.Cells.Item(24, 3).Value2 = 1.051033241513778
$ws.Cells.Item(24, 4).Value2 = 1.056612297002211
$ws.Cells.Item(24, 5).Value2 = 1.054875366993659
$ws.Cells.Item(24, 6).Value2 = 1.065598757727077
$ws.Cells.Item(24, 9).Value2 = 1.051383505770745
$ws.Cells.Item(24, 10).Value2 = 1.057123372144639
$ws.Cells.Item(24, 11).Value2 = 1.059911970139387
$ws.Cells.Item(24, 12).Value2 = 1.058180854935609
$ws.Cells.Item(24, 13).Value2 = 1.068868678493598
$ws.Cells.Item(24, 14).Value2 = 1.022663136333487
$ws.Cells.Item(25, 2).Value2 = 1.02
$ws.Cells.Item(25, 3).Value2 = 1.052983624794947
$ws.Cells.Item(25, 4).Value2 = 1.058158491465744
$ws.Cells.Item(25, 5).Value2 = 1.056736960325974
$ws.Cells.Item(25, 6).Value2 = 1.067386609291185
$ws.Cells.Item(25, 9).Value2 = 1.052068605537553
$ws.Cells.Item(25, 10).Value2 = 1.058474917283641
$ws.Cells.Item(25, 11).Value2 = 1.061142373338211
$ws.Cells.Item(25, 12).Value2 = 1.059725112465238
$ws.Cells.Item(25, 13).Value2 = 1.070343071049085
$ws.Cells.Item(25, 14).Value2 = 1.023135248419324
